$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Add new sources" task text (row 10, column B) to mention Bing/Microsoft ---
$ws.Range("B10").Value = "Add new sources in response to user requests (e.g. specific searching of Notornis, Ostritch and Kukila journal websites, and the abstract database Scopus. Possibly ResearchGate?). Enable Bing search through discussion with MicroSoft again."

# --- New column C: "Must/Should/Could/Will not have for MVP" ---
# Header
$ws.Range("C1").Value = "Must/Should/Could/Will not have for MVP"
$ws.Range("C1").Font.Bold = $true

# MoSCoW priority values for the highlighted (yellow) rows
$ws.Range("C10").Value = "M"
$ws.Range("C13").Value = "S/C"
$ws.Range("C15").Value = "W but good to think through the possibilities here soon"
$ws.Range("C19").Value = "M"
$ws.Range("C22").Value = "S/C"

# Match the yellow highlight already used for A/B on those rows
$ws.Range("C10").Interior.Color = 65535
$ws.Range("C13").Interior.Color = 65535
$ws.Range("C15").Interior.Color = 65535
$ws.Range("C19").Interior.Color = 65535
$ws.Range("C22").Interior.Color = 65535

# --- Column widths: A stays ~same, new column B (task text) becomes wide ---
$ws.Columns.Item(1).ColumnWidth = 12.333333333333334
$ws.Columns.Item(2).ColumnWidth = 97.66666666666667

# --- Restore the cursor/selection position left by the editor ---
[void]$ws.Range("B30").Select()

Write-Output "done"
